{"js": "// The published site rebuilt lom3097.docx and, in doing so, dropped the\n// trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line, the\n// \"\u00a9 2020 ... Creative Commons Attribution\" footer line, and the blank\n// paragraph that separated them from the \"LOM3070: Est\u00e1gio Supervisionado\n// (Requisito)\" line above \u2014 three whole paragraphs removed, leaving the\n// single blank paragraph (and the page-break paragraph) that originally\n// followed the footer untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"LOM3070: ... (Requisito)\") so we only ever\n// touch the specific run of paragraphs that follows it, not any other blank\n// paragraph elsewhere in the document.\nconst anchorText = \"LOM3070: Est\u00e1gio Supervisionado (Requisito)\";\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const targetsInOrder = [\n    \"\", // blank paragraph right after the \"Requisitos\" answer\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n  ];\n\n  const toDelete = [];\n  let cursor = anchorIndex + 1;\n  for (const expected of targetsInOrder) {\n    if (cursor < items.length && items[cursor].text.trim() === expected) {\n      toDelete.push(items[cursor]);\n      cursor++;\n    } else {\n      // Structure didn't match what we expected; bail out rather than\n      // deleting the wrong paragraphs.\n      toDelete.length = 0;\n      break;\n    }\n  }\n\n  for (const paragraph of toDelete) {\n    paragraph.delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# The published site rebuilt lom3097.docx and, in doing so, dropped the\n# trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line, the\n# \"\u00a9 2020 ... Creative Commons Attribution\" footer line, and the blank\n# paragraph that separated them from the \"LOM3070: Est\u00e1gio Supervisionado\n# (Requisito)\" line above \u2014 three whole paragraphs removed, leaving the\n# single blank paragraph (and the page-break paragraph) that originally\n# followed the footer untouched.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"LOM3070: Est\u00e1gio Supervisionado (Requisito)\"\n\n$paras = @($d.Paragraphs)\n\n$anchorIndex = -1\nfor ($i = 0; $i -lt $paras.Count; $i++) {\n    if ($paras[$i].Range.Text.Trim() -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ge 0) {\n    $targetsInOrder = @(\n        \"\",\n        \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n        \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n    )\n\n    $toDelete = @()\n    $cursor = $anchorIndex + 1\n    $ok = $true\n    foreach ($expected in $targetsInOrder) {\n        if ($cursor -lt $paras.Count -and $paras[$cursor].Range.Text.Trim() -eq $expected) {\n            $toDelete += $paras[$cursor]\n            $cursor = $cursor + 1\n        } else {\n            $ok = $false\n            break\n        }\n    }\n\n    if ($ok) {\n        for ($j = $toDelete.Count - 1; $j -ge 0; $j--) {\n            $toDelete[$j].Range.Delete()\n        }\n    }\n}\n"}
